# Auto-generated market-price / profit refresh for Maduin_Profits workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per
# refreshed market-board data. Cells whose recomputed profit no longer
# applies (e.g. the HQ side has no HQ recipe) are cleared instead of zeroed,
# matching how the source sheet represents 'not applicable'.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 476.85715  # was 212
$ws.Range("I4").Value = 476.85715  # was 212
$ws.Range("K4").Value = 476.85715  # was 212
$ws.Range("M4").Value = -362.85715  # was -98
$ws.Range("H6").Value = 203.85715  # was 204.5
$ws.Range("I6").Value = 221.16667  # was 225.4
$ws.Range("K6").Value = 663.50001  # was 676.2
$ws.Range("M6").Value = -551.50001  # was -564.2
$ws.Range("H17").Value = 1375.2  # was 1063
$ws.Range("I17").Value = 1125.3334  # was 1063
$ws.Range("J17").Value = 1750  # was 0
$ws.Range("K17").Value = 3376.0002  # was 3189
$ws.Range("L17").Value = 5250  # was 0
$ws.Range("M17").Value = -3208.0002  # was -3021
$ws.Range("N17").Value = -5586  # was (empty)
$ws.Range("H38").Value = 53.125  # was 54.875
$ws.Range("I38").Value = 53.125  # was 54.875
$ws.Range("K38").Value = 159.375  # was 164.625
$ws.Range("M38").Value = 212.625  # was 207.375
$ws.Range("H39").Value = 1630  # was 454.1111
$ws.Range("I39").Value = 91.22221999999999  # was 121
$ws.Range("J39").Value = 3361.125  # was 720.6
$ws.Range("K39").Value = 273.66666  # was 363
$ws.Range("L39").Value = 10083.375  # was 2161.8
$ws.Range("M39").Value = 22.33334000000002  # was -67
$ws.Range("N39").Value = -10675.375  # was -2753.8
$ws.Range("H86").Value = 3083.1667  # was 10000
$ws.Range("I86").Value = 3083.1667  # was 12000
$ws.Range("J86").Value = 0  # was 7000
$ws.Range("K86").Value = 3083.1667  # was 12000
$ws.Range("L86").Value = 0  # was 7000
$ws.Range("M86").Value = -1960.1667  # was -10877
$ws.Range("N86").ClearContents()  # was -9246
$ws.Range("H88").Value = 4297.8335  # was 4312.2856
$ws.Range("I88").Value = 0  # was 4399
$ws.Range("K88").Value = 0  # was 4399
$ws.Range("M88").ClearContents()  # was -3993
$ws.Range("H89").Value = 3083.1667  # was 10000
$ws.Range("I89").Value = 3083.1667  # was 12000
$ws.Range("J89").Value = 0  # was 7000
$ws.Range("K89").Value = 15415.8335  # was 60000
$ws.Range("L89").Value = 0  # was 35000
$ws.Range("M89").Value = -9799.833500000001  # was -54384
$ws.Range("N89").ClearContents()  # was -46232
$ws.Range("H91").Value = 4297.8335  # was 4312.2856
$ws.Range("I91").Value = 0  # was 4399
$ws.Range("K91").Value = 0  # was 4399
$ws.Range("M91").ClearContents()  # was -2995
$ws.Range("H105").Value = 0  # was 11200
$ws.Range("I105").Value = 0  # was 11200
$ws.Range("K105").Value = 0  # was 11200
$ws.Range("M105").ClearContents()  # was -7706
$ws.Range("H106").Value = 0  # was 1005
$ws.Range("I106").Value = 0  # was 1005
$ws.Range("K106").Value = 0  # was 1005
$ws.Range("M106").ClearContents()  # was -374
$ws.Range("H137").Value = 2199.4  # was 2592.6428
$ws.Range("I137").Value = 699.1429000000001  # was 716.1667
$ws.Range("J137").Value = 3512.125  # was 4000
$ws.Range("K137").Value = 2097.4287  # was 2148.5001
$ws.Range("L137").Value = 10536.375  # was 12000
$ws.Range("M137").Value = 452.5712999999996  # was 401.4998999999998
$ws.Range("N137").Value = -15636.375  # was -17100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2333  # was 2499.5
$ws.Range("J74").Value = 2000  # was 0
$ws.Range("L74").Value = 2000  # was 0
$ws.Range("N74").Value = -3748  # was (empty)
$ws.Range("H77").Value = 2333  # was 2499.5
$ws.Range("J77").Value = 2000  # was 0
$ws.Range("L77").Value = 10000  # was 0
$ws.Range("N77").Value = -18736  # was (empty)
$ws.Range("H94").Value = 22999.5  # was 0
$ws.Range("J94").Value = 22999.5  # was 0
$ws.Range("L94").Value = 22999.5  # was 0
$ws.Range("N94").Value = -24801.5  # was (empty)
$ws.Range("H95").Value = 0  # was 7500
$ws.Range("J95").Value = 0  # was 7500
$ws.Range("L95").Value = 0  # was 7500
$ws.Range("N95").ClearContents()  # was -12992
$ws.Range("H98").Value = 33397.8  # was 39331.668
$ws.Range("J98").Value = 33397.8  # was 39331.668
$ws.Range("L98").Value = 33397.8  # was 39331.668
$ws.Range("N98").Value = -39387.8  # was -45321.668
$ws.Range("H132").Value = 717.4286  # was 735.1667
$ws.Range("I132").Value = 717.4286  # was 735.1667
$ws.Range("K132").Value = 2152.2858  # was 2205.5001
$ws.Range("M132").Value = 377.7142000000003  # was 324.4998999999998
$ws.Range("H140").Value = 45000  # was 46666.332
$ws.Range("J140").Value = 60000  # was 54999.5
$ws.Range("L140").Value = 60000  # was 54999.5
$ws.Range("N140").Value = -70360  # was -65359.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 11833.333  # was 11875
$ws.Range("J100").Value = 11833.333  # was 11875
$ws.Range("L100").Value = 11833.333  # was 11875
$ws.Range("N100").Value = -13997.333  # was -14039

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 12899  # was 16821.5
$ws.Range("J28").Value = 12899  # was 16821.5
$ws.Range("L28").Value = 12899  # was 16821.5
$ws.Range("N28").Value = -13389  # was -17311.5
$ws.Range("H31").Value = 2775.1333  # was 2969.5833
$ws.Range("I31").Value = 2708.818  # was 2867.111
$ws.Range("J31").Value = 2957.5  # was 3277
$ws.Range("K31").Value = 2708.818  # was 2867.111
$ws.Range("L31").Value = 2957.5  # was 3277
$ws.Range("M31").Value = -2413.818  # was -2572.111
$ws.Range("N31").Value = -3547.5  # was -3867
$ws.Range("H32").Value = 3427.4285  # was 3848.75
$ws.Range("I32").Value = 3427.4285  # was 3848.75
$ws.Range("K32").Value = 3427.4285  # was 3848.75
$ws.Range("M32").Value = -3111.4285  # was -3532.75
$ws.Range("H34").Value = 2775.1333  # was 2969.5833
$ws.Range("I34").Value = 2708.818  # was 2867.111
$ws.Range("J34").Value = 2957.5  # was 3277
$ws.Range("K34").Value = 2708.818  # was 2867.111
$ws.Range("L34").Value = 2957.5  # was 3277
$ws.Range("M34").Value = -2506.818  # was -2665.111
$ws.Range("N34").Value = -3361.5  # was -3681
$ws.Range("H58").Value = 896.5  # was 898
$ws.Range("I58").Value = 896.5  # was 897
$ws.Range("J58").Value = 0  # was 899
$ws.Range("K58").Value = 896.5  # was 897
$ws.Range("L58").Value = 0  # was 899
$ws.Range("M58").Value = -693.5  # was -694
$ws.Range("N58").ClearContents()  # was -1305
$ws.Range("H134").Value = 4051.75  # was 3940.8333
$ws.Range("I134").Value = 3552.3  # was 3419.2
$ws.Range("K134").Value = 10656.9  # was 10257.6
$ws.Range("M134").Value = -8121.900000000001  # was -7722.599999999999
$ws.Range("H136").Value = 896.5  # was 898
$ws.Range("I136").Value = 896.5  # was 897
$ws.Range("J136").Value = 0  # was 899
$ws.Range("K136").Value = 2689.5  # was 2691
$ws.Range("L136").Value = 0  # was 2697
$ws.Range("M136").Value = -139.5  # was -141
$ws.Range("N136").ClearContents()  # was -7797

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 400.83334  # was 529.2857
$ws.Range("I5").Value = 281.4  # was 451.16666
$ws.Range("K5").Value = 844.1999999999999  # was 1353.49998
$ws.Range("M5").Value = -732.1999999999999  # was -1241.49998
$ws.Range("H135").Value = 400.83334  # was 529.2857
$ws.Range("I135").Value = 281.4  # was 451.16666
$ws.Range("K135").Value = 2532.6  # was 4060.49994
$ws.Range("M135").Value = 2.400000000000091  # was -1525.49994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 18333.334  # was 0
$ws.Range("J39").Value = 18333.334  # was 0
$ws.Range("L39").Value = 18333.334  # was 0
$ws.Range("N39").Value = -19397.334  # was (empty)
$ws.Range("H92").Value = 12127.4  # was 11106.167
$ws.Range("J92").Value = 12127.4  # was 11106.167
$ws.Range("L92").Value = 12127.4  # was 11106.167
$ws.Range("N92").Value = -15871.4  # was -14850.167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 0  # was 1950
$ws.Range("I32").Value = 0  # was 1950
$ws.Range("K32").Value = 0  # was 1950
$ws.Range("M32").ClearContents()  # was -1633
$ws.Range("H94").Value = 25165  # was 25822.5
$ws.Range("J94").Value = 25165  # was 25822.5
$ws.Range("L94").Value = 25165  # was 25822.5
$ws.Range("N94").Value = -26517  # was -27174.5
$ws.Range("H132").Value = 5600  # was 5318
$ws.Range("I132").Value = 5600  # was 5619.9
$ws.Range("J132").Value = 0  # was 2299
$ws.Range("K132").Value = 16800  # was 16859.7
$ws.Range("L132").Value = 0  # was 6897
$ws.Range("M132").Value = -14270  # was -14329.7
$ws.Range("N132").ClearContents()  # was -11957

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2962.3125  # was 2846.8823
$ws.Range("I132").Value = 2924.3572  # was 2796.0667
$ws.Range("K132").Value = 8773.071599999999  # was 8388.2001
$ws.Range("M132").Value = -6243.071599999999  # was -5858.2001
$ws.Range("H136").Value = 3564.4167  # was 3374.8462
$ws.Range("I136").Value = 2524.818  # was 2406.0833
$ws.Range("K136").Value = 7574.454000000001  # was 7218.249899999999
$ws.Range("M136").Value = -5024.454000000001  # was -4668.249899999999

